$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3950
$ws.Range("I51").Value = 6000
$ws.Range("J51").Value = 1900
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 1900
$ws.Range("M51").Value = -5516
$ws.Range("N51").Value = -2868

$ws.Range("H96").Value = 1010.8571
$ws.Range("I96").Value = 1131.5
$ws.Range("J96").Value = 850
$ws.Range("K96").Value = 3394.5
$ws.Range("L96").Value = 2550
$ws.Range("M96").Value = -2021.5
$ws.Range("N96").Value = -5296

$ws.Range("H132").Value = 2760.75
$ws.Range("I132").Value = 2699.4062
$ws.Range("K132").Value = 8098.2186
$ws.Range("M132").Value = -5568.2186

$ws.Range("H138").Value = 1682.3636
$ws.Range("I138").Value = 937.5333000000001
$ws.Range("J138").Value = 2729.7812
$ws.Range("K138").Value = 2812.5999
$ws.Range("L138").Value = 8189.3436
$ws.Range("M138").Value = 2327.4001
$ws.Range("N138").Value = -18469.3436

$ws.Range("H141").Value = 2643.9285
$ws.Range("I141").Value = 1949.375
$ws.Range("J141").Value = 3570
$ws.Range("K141").Value = 5848.125
$ws.Range("L141").Value = 10710
$ws.Range("M141").Value = -668.125
$ws.Range("N141").Value = -21070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27597.879
$ws.Range("I32").Value = 28789.447
$ws.Range("K32").Value = 28789.447
$ws.Range("M32").Value = -28502.447

$ws.Range("H45").Value = 3014.3333
$ws.Range("I45").Value = 4250.1665
$ws.Range("J45").Value = 2190.4443
$ws.Range("K45").Value = 4250.1665
$ws.Range("L45").Value = 2190.4443
$ws.Range("M45").Value = -3873.1665
$ws.Range("N45").Value = -2944.4443

$ws.Range("H63").Value = 1550
$ws.Range("I63").Value = 1550
$ws.Range("K63").Value = 1550
$ws.Range("M63").Value = -864

$ws.Range("H66").Value = 1550
$ws.Range("I66").Value = 1550
$ws.Range("K66").Value = 7750
$ws.Range("M66").Value = -4318

$ws.Range("H102").Value = 1731.2354
$ws.Range("I102").Value = 1436.8889
$ws.Range("J102").Value = 2062.375
$ws.Range("K102").Value = 1436.8889
$ws.Range("L102").Value = 2062.375
$ws.Range("M102").Value = 185.1111000000001
$ws.Range("N102").Value = -5306.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2479.3333
$ws.Range("I20").Value = 2911.842
$ws.Range("J20").Value = 1452.125
$ws.Range("K20").Value = 2911.842
$ws.Range("L20").Value = 1452.125
$ws.Range("M20").Value = -2664.842
$ws.Range("N20").Value = -1946.125

$ws.Range("H82").Value = 11897.909
$ws.Range("I82").Value = 6534.2
$ws.Range("J82").Value = 65535
$ws.Range("K82").Value = 6534.2
$ws.Range("L82").Value = 65535
$ws.Range("M82").Value = -6151.2
$ws.Range("N82").Value = -66301

$ws.Range("H85").Value = 11897.909
$ws.Range("I85").Value = 6534.2
$ws.Range("J85").Value = 65535
$ws.Range("K85").Value = 6534.2
$ws.Range("L85").Value = 65535
$ws.Range("M85").Value = -5208.2
$ws.Range("N85").Value = -68187

$ws.Range("H99").Value = 1513.7333
$ws.Range("I99").Value = 1233
$ws.Range("K99").Value = 1233
$ws.Range("M99").Value = 265

$ws.Range("H107").Value = 1070
$ws.Range("I107").Value = 1049.5
$ws.Range("K107").Value = 1049.5
$ws.Range("M107").Value = 870.5

$ws.Range("H134").Value = 43092.81
$ws.Range("I134").Value = 48344.043
$ws.Range("J134").Value = 2833.3333
$ws.Range("K134").Value = 145032.129
$ws.Range("L134").Value = 8499.999899999999
$ws.Range("M134").Value = -142497.129
$ws.Range("N134").Value = -13569.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9741.737999999999
$ws.Range("I31").Value = 12138.034
$ws.Range("J31").Value = 4396.154
$ws.Range("K31").Value = 12138.034
$ws.Range("L31").Value = 4396.154
$ws.Range("M31").Value = -11843.034
$ws.Range("N31").Value = -4986.154

$ws.Range("H34").Value = 9741.737999999999
$ws.Range("I34").Value = 12138.034
$ws.Range("J34").Value = 4396.154
$ws.Range("K34").Value = 12138.034
$ws.Range("L34").Value = 4396.154
$ws.Range("M34").Value = -11936.034
$ws.Range("N34").Value = -4800.154

$ws.Range("H39").Value = 9275.5
$ws.Range("J39").Value = 17500
$ws.Range("L39").Value = 17500
$ws.Range("N39").Value = -18282

$ws.Range("H49").Value = 9275.5
$ws.Range("J49").Value = 17500
$ws.Range("L49").Value = 17500
$ws.Range("N49").Value = -17864

$ws.Range("H58").Value = 14551.189
$ws.Range("I58").Value = 1051.963
$ws.Range("K58").Value = 1051.963
$ws.Range("M58").Value = -848.963

$ws.Range("H74").Value = 34693.832
$ws.Range("J74").Value = 34693.832
$ws.Range("L74").Value = 34693.832
$ws.Range("N74").Value = -36441.832

$ws.Range("H77").Value = 34693.832
$ws.Range("J77").Value = 34693.832
$ws.Range("L77").Value = 104081.496
$ws.Range("N77").Value = -112817.496

$ws.Range("H88").Value = 36136.6
$ws.Range("J88").Value = 36136.6
$ws.Range("L88").Value = 36136.6
$ws.Range("N88").Value = -36948.6

$ws.Range("H91").Value = 36136.6
$ws.Range("J91").Value = 36136.6
$ws.Range("L91").Value = 36136.6
$ws.Range("N91").Value = -38944.6

$ws.Range("H136").Value = 14551.189
$ws.Range("I136").Value = 1051.963
$ws.Range("K136").Value = 3155.889
$ws.Range("M136").Value = -605.8890000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1794.3334
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 1912
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 5736
$ws.Range("M46").Value = -1409
$ws.Range("N46").Value = -5918

$ws.Range("H122").Value = 773.2308
$ws.Range("J122").Value = 933.55554
$ws.Range("L122").Value = 8401.99986
$ws.Range("N122").Value = -13301.99986

$ws.Range("H131").Value = 738.59
$ws.Range("I131").Value = 533.1667
$ws.Range("J131").Value = 751.70215
$ws.Range("K131").Value = 1599.5001
$ws.Range("L131").Value = 2255.10645
$ws.Range("M131").Value = 3440.4999
$ws.Range("N131").Value = -12335.10645

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 53890.133
$ws.Range("I132").Value = 51879.285
$ws.Range("J132").Value = 58582.11
$ws.Range("K132").Value = 155637.855
$ws.Range("L132").Value = 175746.33
$ws.Range("M132").Value = -153107.855
$ws.Range("N132").Value = -180806.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 16333.333
$ws.Range("I41").Value = 15000
$ws.Range("J41").Value = 17000
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 17000
$ws.Range("M41").Value = -14562
$ws.Range("N41").Value = -17876

$ws.Range("H46").Value = 1177.6666
$ws.Range("I46").Value = 874.75
$ws.Range("J46").Value = 1420
$ws.Range("K46").Value = 874.75
$ws.Range("L46").Value = 1420
$ws.Range("M46").Value = -686.75
$ws.Range("N46").Value = -1796

$ws.Range("H61").Value = 4318.9375
$ws.Range("I61").Value = 3110
$ws.Range("J61").Value = 6333.8335
$ws.Range("K61").Value = 3110
$ws.Range("L61").Value = 6333.8335
$ws.Range("M61").Value = -2908
$ws.Range("N61").Value = -6737.8335

$ws.Range("H113").Value = 4318.9375
$ws.Range("I113").Value = 3110
$ws.Range("J113").Value = 6333.8335
$ws.Range("K113").Value = 3110
$ws.Range("L113").Value = 6333.8335
$ws.Range("M113").Value = -940
$ws.Range("N113").Value = -10673.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 795
$ws.Range("I21").Value = 452.77777
$ws.Range("J21").Value = 1565
$ws.Range("K21").Value = 452.77777
$ws.Range("L21").Value = 1565
$ws.Range("M21").Value = -217.77777
$ws.Range("N21").Value = -2035

$ws.Range("H35").Value = 795
$ws.Range("I35").Value = 452.77777
$ws.Range("J35").Value = 1565
$ws.Range("K35").Value = 452.77777
$ws.Range("L35").Value = 1565
$ws.Range("M35").Value = -162.77777
$ws.Range("N35").Value = -2145

$ws.Range("H81").Value = 125001580
$ws.Range("I81").Value = 2260
$ws.Range("J81").Value = 333333800
$ws.Range("K81").Value = 4520
$ws.Range("L81").Value = 666667600
$ws.Range("M81").Value = -3459
$ws.Range("N81").Value = -666669722

$ws.Range("H84").Value = 125001580
$ws.Range("I84").Value = 2260
$ws.Range("J84").Value = 333333800
$ws.Range("K84").Value = 22600
$ws.Range("L84").Value = 3333338000
$ws.Range("M84").Value = -17296
$ws.Range("N84").Value = -3333348608

$ws.Range("H113").Value = 1590616.5
$ws.Range("I113").Value = 1105.4546
$ws.Range("J113").Value = 4504720
$ws.Range("K113").Value = 3316.3638
$ws.Range("L113").Value = 13514160
$ws.Range("M113").Value = -1146.3638
$ws.Range("N113").Value = -13518500
